# Weekly fruit/vegetable update: insert a new record row for Caigua
# (Agricola del Norte S.A. de Arica) at row 107, pushing the existing
# rows 107-111 down to 108-112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107 (shifts old rows 107-111 down to 108-112,
# carrying their formatting/styles along).
$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with the new weekly record.
$ws.Range("A107").Value = 1
$ws.Range("B107").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C107").Value = "Arica y Parinacota"
$ws.Range("D107").Value = 44746
$ws.Range("E107").Value = 15
$ws.Range("F107").Value = 100112036
$ws.Range("G107").Value = "Caigua"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 120
$ws.Range("K107").Value = 7000
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = 7500
$ws.Range("N107").Value = "$/caja 20 kilos"
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 375
$ws.Range("Q107").Value = 20
$ws.Range("R107").Value = "Hortaliza"
